$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The underlying data source got reloaded and a handful of rows now carry
# different (but still duplicate-grouped) values - course name/abbreviation
# for "Algoritmy a datove struktury" / "Teorie formalnich jazyku" swapped
# between rows 6 and 7, and several abbreviation columns got re-ordered
# for their duplicate-name groups.

$ws.Range("A6").Value = "Algoritmy a datové struktury"
$ws.Range("B6").Value = "KDSA"
$ws.Range("D6").Value = "''doc. RNDr. Karel Oliva Dr."

$ws.Range("A7").Value = "Teorie formálních jazyků"
$ws.Range("B7").Value = "TFL"
$ws.Range("D7").Value = "''doc. RNDr. Karel Oliva Dr.'"

$ws.Range("B8").Value = "DSA"

$ws.Range("B9").Value = "KSYS"
$ws.Range("B10").Value = "SYS"

$ws.Range("B12").Value = "KBIG"
$ws.Range("B13").Value = "BIG"

$ws.Range("B18").Value = "KANE"
$ws.Range("B19").Value = "KAEL"
